$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15: Log / GET / study with Args "studyId=1234"
$ws.Range("D15").Value = "studyId=1234"
$ws.Range("E15:I15").Clear()
$ws.Range("F15").Formula = "=CONCATENATE(""http://wlux.uw.edu/data/"",A15,"".php"", IF(D15<>"""",""?"",""""),D15)"
$ws.Range("G15").Formula = "=CONCATENATE(""http://localhost/data/"",A15,"".php"", IF(D15<>"""",""?"",""""),D15)"

# Row 16: Log / GET / session with Args "sessionId=1383605381"
$ws.Range("D16").Value = "sessionId=1383605381"
$ws.Range("E16:I16").Clear()
$ws.Range("F16").Formula = "=CONCATENATE(""http://wlux.uw.edu/data/"",A16,"".php"", IF(D16<>"""",""?"",""""),D16)"
$ws.Range("G16").Formula = "=CONCATENATE(""http://localhost/data/"",A16,"".php"", IF(D16<>"""",""?"",""""),D16)"

# Row 17: Log / GET / task with Args "sessionId=1383605381&taskId=1"
$ws.Range("D17").Value = "sessionId=1383605381&taskId=1"
$ws.Range("E17:I17").Clear()
$ws.Range("F17").Formula = "=CONCATENATE(""http://wlux.uw.edu/data/"",A17,"".php"", IF(D17<>"""",""?"",""""),D17)"
$ws.Range("G17").Formula = "=CONCATENATE(""http://localhost/data/"",A17,"".php"", IF(D17<>"""",""?"",""""),D17)"

# view changes
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F17").Select()
